$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "58.625.29"
Set-TextValue "E2" "  -2.66%  "
Set-TextValue "D3" "2.715.72"
Set-TextValue "E3" "  -6.43%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "503.77"
Set-TextValue "E5" "  -4.52%  "
Set-TextValue "D6" "139.20"
Set-TextValue "E6" "  -1.64%  "
Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.29%  "
Set-TextValue "D8" "0.529"
Set-TextValue "E8" "  -3.56%  "
Set-TextValue "D9" "2.724.78"
Set-TextValue "E9" "  -6.10%  "
Set-TextValue "D10" "6.11"
Set-TextValue "E10" "  +4.63%  "
Set-TextValue "E11" "  -3.61%  "
Set-TextValue "D12" "0.345"
Set-TextValue "E12" "  -1.90%  "
Set-TextValue "E13" "  +1.17%  "
Set-TextValue "D14" "3.183.21"
Set-TextValue "E14" "  -6.55%  "
Set-TextValue "D15" "58.546.79"
Set-TextValue "E15" "  -3.06%  "
Set-TextValue "D16" "21.46"
Set-TextValue "E16" "  -5.04%  "
Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "2.711.32"
Set-TextValue "E17" "  -6.56%  "
Set-TextValue "B18" "ShibaInu"
Set-TextValue "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D18" "0.0000134"
Set-TextValue "E18" "  -3.71%  "
Set-TextValue "D19" "4.68"
Set-TextValue "E19" "  -4.40%  "
Set-TextValue "D20" "10.87"
Set-TextValue "E20" "  -5.11%  "
Set-TextValue "D21" "341.66"
Set-TextValue "E21" "  -4.93%  "
Set-TextValue "D22" "6.18"
Set-TextValue "E22" "  -5.73%  "
Set-TextValue "E23" "  -0.03%  "
Set-TextValue "D24" "5.61"
Set-TextValue "E24" "  -0.60%  "
Set-TextValue "D25" "62.56"
Set-TextValue "E25" "  -1.23%  "
Set-TextValue "D26" "0.422"
Set-TextValue "E26" "  -5.33%  "
Set-TextValue "D27" "0.170"
Set-TextValue "E27" "  -2.53%  "
Set-TextValue "E28" "  -0.60%  "
Set-TextValue "D29" "0.0₃0822"
Set-TextValue "E29" "  -2.65%  "
Set-TextValue "D30" "7.38"
Set-TextValue "E30" "  -3.70%  "
Set-TextValue "D31" "0.998"
Set-TextValue "E31" "  -0.15%  "
Set-TextValue "D32" "1.60"
Set-TextValue "E32" "  -3.72%  "
Set-TextValue "D33" "19.00"
Set-TextValue "E33" "  -2.49%  "
Set-TextValue "D34" "148.48"
Set-TextValue "E34" "  -0.87%  "
Set-TextValue "D35" "4.15"
Set-TextValue "E35" "  -2.84%  "
Set-TextValue "D36" "5.32"
Set-TextValue "E36" "  -3.23%  "
Set-TextValue "D37" "0.929"
Set-TextValue "E37" "  -6.06%  "
Set-TextValue "E38" "  -5.12%  "
Set-TextValue "D39" "35.99"
Set-TextValue "E39" "  -4.19%  "
Set-TextValue "E40" "  -5.28%  "
Set-TextValue "D41" "2.168.97"
Set-TextValue "E41" "  -6.97%  "
Set-TextValue "D42" "3.49"
Set-TextValue "E42" "  -4.00%  "
Set-TextValue "D43" "0.996"
Set-TextValue "E43" "  -0.32%  "
Set-TextValue "E44" "  -2.66%  "
Set-TextValue "D45" "0.598"
Set-TextValue "E45" "  -6.68%  "
Set-TextValue "D46" "18.78"
Set-TextValue "E46" "  -9.09%  "
Set-TextValue "D47" "10.35"
Set-TextValue "E47" "  -0.01%  "
Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "4.72"
Set-TextValue "E48" "  -2.06%  "
Set-TextValue "B49" "VeChain"
Set-TextValue "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0225"
Set-TextValue "E49" "  -2.89%  "
Set-TextValue "E50" "  -4.47%  "
Set-TextValue "D51" "17.73"
Set-TextValue "E51" "  -2.11%  "
